$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7: Logout / CLICK / welcome / id ---
$ws.Range("A7:E7").Borders.LineStyle = 1
$ws.Range("A7").Value = "Logout"
$ws.Range("B7").Value = "CLICK"
$ws.Range("C7").Value = "welcome"
$ws.Range("D7").Value = "id"

# --- New row 8: (blank) / CLICK / logout / css ---
$ws.Range("A8:E8").Borders.LineStyle = 1
$ws.Range("B8").Value = "CLICK"
$ws.Range("C8").Value = "logout"
$ws.Range("D8").Value = "css"

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
